$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DAMSLTag (I) and DialogAct (J) columns per re-run of SGNN dialog-act annotation.
$ws.Range("I2").Value = 'b'
$ws.Range("J2").Value = 'Acknowledge (Backchannel)'
$ws.Range("I36").Value = 'sd'
$ws.Range("J36").Value = 'Statement-non-opinion'
$ws.Range("I39").Value = '%'
$ws.Range("J39").Value = 'Uninterpretable'
$ws.Range("I41").Value = 'sv'
$ws.Range("J41").Value = 'Statement-opinion'
$ws.Range("I49").Value = 'b'
$ws.Range("J49").Value = 'Acknowledge (Backchannel)'
$ws.Range("I52").Value = 'aa'
$ws.Range("J52").Value = 'Agree/Accept'
$ws.Range("I56").Value = 'b'
$ws.Range("J56").Value = 'Acknowledge (Backchannel)'
$ws.Range("I81").Value = 'aa'
$ws.Range("J81").Value = 'Agree/Accept'
$ws.Range("I96").Value = 'sd'
$ws.Range("J96").Value = 'Statement-non-opinion'
$ws.Range("I108").Value = 'sv'
$ws.Range("J108").Value = 'Statement-opinion'
$ws.Range("I133").Value = 'sv'
$ws.Range("J133").Value = 'Statement-opinion'
$ws.Range("I139").Value = 'sd'
$ws.Range("J139").Value = 'Statement-non-opinion'
$ws.Range("I140").Value = 'sd'
$ws.Range("J140").Value = 'Statement-non-opinion'
$ws.Range("I146").Value = 'sv'
$ws.Range("J146").Value = 'Statement-opinion'
$ws.Range("I152").Value = 'sd'
$ws.Range("J152").Value = 'Statement-non-opinion'
$ws.Range("I161").Value = 'sv'
$ws.Range("J161").Value = 'Statement-opinion'
$ws.Range("I168").Value = 'sv'
$ws.Range("J168").Value = 'Statement-opinion'
$ws.Range("I171").Value = 'sv'
$ws.Range("J171").Value = 'Statement-opinion'
$ws.Range("I172").Value = 'ba'
$ws.Range("J172").Value = 'Appreciation'
$ws.Range("I180").Value = 'sd'
$ws.Range("J180").Value = 'Statement-non-opinion'
$ws.Range("I191").Value = 'sd'
$ws.Range("J191").Value = 'Statement-non-opinion'
$ws.Range("I192").Value = 'aa'
$ws.Range("J192").Value = 'Agree/Accept'
$ws.Range("I193").Value = 'b'
$ws.Range("J193").Value = 'Acknowledge (Backchannel)'
$ws.Range("I198").Value = 'sd'
$ws.Range("J198").Value = 'Statement-non-opinion'
$ws.Range("I227").Value = 'sd'
$ws.Range("J227").Value = 'Statement-non-opinion'
$ws.Range("I262").Value = 'ba'
$ws.Range("J262").Value = 'Appreciation'
$ws.Range("I286").Value = 'b'
$ws.Range("J286").Value = 'Acknowledge (Backchannel)'
$ws.Range("I330").Value = 'sd'
$ws.Range("J330").Value = 'Statement-non-opinion'
$ws.Range("I333").Value = 'sd'
$ws.Range("J333").Value = 'Statement-non-opinion'
$ws.Range("I344").Value = 'sd'
$ws.Range("J344").Value = 'Statement-non-opinion'
$ws.Range("I364").Value = 'sv'
$ws.Range("J364").Value = 'Statement-opinion'
$ws.Range("I381").Value = 'sd'
$ws.Range("J381").Value = 'Statement-non-opinion'
$ws.Range("I402").Value = 'sd'
$ws.Range("J402").Value = 'Statement-non-opinion'
$ws.Range("I426").Value = 'sd'
$ws.Range("J426").Value = 'Statement-non-opinion'
$ws.Range("I454").Value = 'sv'
$ws.Range("J454").Value = 'Statement-opinion'
$ws.Range("I466").Value = 'sv'
$ws.Range("J466").Value = 'Statement-opinion'
$ws.Range("I483").Value = 'sd'
$ws.Range("J483").Value = 'Statement-non-opinion'
$ws.Range("I488").Value = 'b'
$ws.Range("J488").Value = 'Acknowledge (Backchannel)'
$ws.Range("I495").Value = 'sd'
$ws.Range("J495").Value = 'Statement-non-opinion'
$ws.Range("I504").Value = 'sd'
$ws.Range("J504").Value = 'Statement-non-opinion'
$ws.Range("I526").Value = 'sd'
$ws.Range("J526").Value = 'Statement-non-opinion'
$ws.Range("I529").Value = 'sd'
$ws.Range("J529").Value = 'Statement-non-opinion'
